$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New station-load data (columns B..I, rows 2..14)
$data = @(
    @(19192,5435,8520,0,5088,0,15274,14868),
    @(18619,5008,8260,0,4996,0,14694,14177),
    @(17969,4808,8031,0,4678,0,14159,13807),
    @(17223,4700,7899,0,4771,0,13738,13447),
    @(17898,4889,8098,0,4917,0,13890,14009),
    @(18664,5231,8714,0,5123,0,14547,14642),
    @(18237,4552,8627,0,5517,0,14033,14593),
    @(21576,4438,9964,0,9197,0,15693,16434),
    @(26902,5153,12410,0,12608,0,18857,18423),
    @(27846,5591,17153,0,13638,0,0,20468),
    @(28721,5830,18641,0,13864,0,17260,21035),
    @(28957,5853,18818,0,14028,0,17640,21320),
    @(27430,5862,19220,0,14041,0,17840,0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j  # column B is index 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
